$d = $word.ActiveDocument

# Insert a brand-new paragraph before the current first paragraph.
# Using InsertParagraphBefore() on the existing first paragraph's range
# duplicates that paragraph's formatting (pStyle "Heading2" plus the nil
# pBdr borders), matching the target markup exactly.
$firstPara = $d.Paragraphs.Item(1)
$firstPara.Range.InsertParagraphBefore()

# Fill in the newly created (now first) paragraph's text.
$newPara = $d.Paragraphs.Item(1)
$newPara.Range.Text = "This file has been modified by Isaac Hanna for the purposes of the in-class 471 assignment 4/22/2022"
